$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Structural changes on Sheet1 ---
# Old layout: A=Alarm B=IType C=Dependence D=bw E=Delta F=Threshold G=nu H=w I=k ... K=Threshold(=F*1e6 formula)
# New layout: A=alarmGen B=infPeriod C=smoothWindow D=beta(new) E=delta F=H G=nu H=x0 I=k

# 1) Remove the old "Dependence" column (C); bw/Delta/Threshold/nu/w/k shift left into C..H
$ws1.Columns.Item(3).Delete()
# 2) Remove the old calculated helper column (was K, now at J after the first delete)
$ws1.Columns.Item(10).Delete()
# 3) Insert a fresh column D to hold the new "beta" parameter
$ws1.Columns.Item(4).Insert()

# --- Header row (row 1) ---
$ws1.Range("A1").Value = "alarmGen"
$ws1.Range("B1").Value = "infPeriod"
$ws1.Range("C1").Value = "smoothWindow"
$ws1.Range("D1").Value = "beta"
$ws1.Range("E1").Value = "delta"
$ws1.Range("F1").Value = "H"
$ws1.Range("G1").Value = "nu"
$ws1.Range("H1").Value = "x0"
$ws1.Range("I1").Value = "k"

# --- Re-label alarm type / dependence-type text values (lower-cased now) ---
$ws1.Range("A2").Value = "thresh"
$ws1.Range("A3").Value = "thresh"
$ws1.Range("A4").Value = "thresh"
$ws1.Range("A5").Value = "thresh"
$ws1.Range("A6").Value = "hill"
$ws1.Range("A7").Value = "hill"
$ws1.Range("A8").Value = "hill"
$ws1.Range("A9").Value = "hill"
$ws1.Range("A10").Value = "power"
$ws1.Range("A11").Value = "power"
$ws1.Range("A12").Value = "power"
$ws1.Range("A13").Value = "power"

$ws1.Range("B2").Value = "fixed"
$ws1.Range("B3").Value = "fixed"
$ws1.Range("B4").Value = "exp"
$ws1.Range("B5").Value = "exp"
$ws1.Range("B6").Value = "fixed"
$ws1.Range("B7").Value = "fixed"
$ws1.Range("B8").Value = "exp"
$ws1.Range("B9").Value = "exp"
$ws1.Range("B10").Value = "fixed"
$ws1.Range("B11").Value = "fixed"
$ws1.Range("B12").Value = "exp"
$ws1.Range("B13").Value = "exp"

# --- Fill the new "beta" column (D) with its constant value for every data row ---
$ws1.Range("D2").Value = 0.36
$ws1.Range("D3").Value = 0.36
$ws1.Range("D4").Value = 0.36
$ws1.Range("D5").Value = 0.36
$ws1.Range("D6").Value = 0.36
$ws1.Range("D7").Value = 0.36
$ws1.Range("D8").Value = 0.36
$ws1.Range("D9").Value = 0.36
$ws1.Range("D10").Value = 0.36
$ws1.Range("D11").Value = 0.36
$ws1.Range("D12").Value = 0.36
$ws1.Range("D13").Value = 0.36

# --- Style touch-ups so the new column D matches its row's existing number format ---
$ws1.Range("D2,D3").Style = $ws1.Range("C2").Style
$ws1.Range("D6,D7,D10,D11").Style = $ws1.Range("E6").Style
$ws1.Range("D4,D5,D8,D9,D12,D13").Style = $ws1.Range("E4").Style

$ws1.Range("A4,A5").Font.ColorIndex = $ws1.Range("A2").Font.ColorIndex
$ws1.Range("A8,A9,A12,A13").Interior.Color = $ws1.Range("A4").Interior.Color
$ws1.Range("B6,B7,B10,B11").Font.Bold = $ws1.Range("B2").Font.Bold

# --- Move selection like the authored workbook ---
$ws1.Activate()
$ws1.Range("H8").Select()

# --- Workbook-level view metadata tweak recorded in the commit ---
$excel.ActiveWindow.WindowState = -4143
$excel.Width = 1280
